$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Alshimaa Atef, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Shimaa Ahmad Mekki, Dr. Heba Mahmoud Ali"
$ws.Range("G3").Value = "Dr. Gehan Adel, Administrator, Dr. Manar Montaser, Dr. Alshimaa Atef"
$ws.Range("G4").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Majorelle Magdy, Dr. Shimaa Ahmad Mekki, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Asmaa Reda, Dr. Hanan Ragab"
$ws.Range("G5").Value = "Dr. Abeer Ragab, Dr. Menna tu'Alllah Mohammad"
$ws.Range("G6").Value = "Dr. Sara Nabil, Dr. Safa Hany"
$ws.Range("G7").Value = "Dr. Amal Awwad, Dr. Safa Hany"
$ws.Range("G9").Value = "Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Eman M. Abo-Sakaya, Dr. Marina Youhanna"
$ws.Range("G10").Value = "Dr. Amira Ibrahim, Dr. Basma Hamed"
$ws.Range("G17").Value = "Dr. Marian Samir, Dr. Enas Omran, Dr. Walaa Ghanima"
$ws.Range("G18").Value = "Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida"
$ws.Range("G19").Value = "Dr. Marina Sorial, Dr. Wafaa Ebida, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Eman Samir Gabry"
$ws.Range("G20").Value = "Dr. Marina Sorial, Dr. Nardine, Dr. Marina Atef, Dr. Neveen Nashaat, Dr. Remon, Dr. Yasmin, Dr. Monica"
$ws.Range("G21").Value = "Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Alshimaa Atef, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Shimaa Ahmad Mekki, Dr. Heba Mahmoud Ali"
$ws.Range("G22").Value = "Dr. Gehan Adel, Administrator, Dr. Manar Montaser, Dr. Alshimaa Atef"
$ws.Range("G23").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Majorelle Magdy, Dr. Shimaa Ahmad Mekki, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Asmaa Reda, Dr. Hanan Ragab"
$ws.Range("G24").Value = "Dr. Fatma Elhady, Dr. Nada Mohammad, Dr. Abeer Ragab, Dr. Lamiaa Ossama, Dr. Amera Ahmad Saad"
$ws.Range("G25").Value = "Dr. Yasmin Tarek, Dr. Nourhan Mohammad"
$ws.Range("G26").Value = "Dr. Amal Awwad, Dr. Safa Hany"
$ws.Range("G28").Value = "Dr. Nourhan Osama, Dr. Marwa Mustafa, Dr. Yasmeena Fattoh, Dr. Sarah Abdelmohsen, Dr. Madeha Saeed, Dr. Arwa Al-Sayed, Dr. Basma Hamed, Dr. Dina Adel, Dr. Eman M. Abo-Sakaya, Dr. Esraa Mostafa"
$ws.Range("G29").Value = "Dr. Amira Ibrahim, Dr. Yasmeena Fattoh, Dr. Esraa Mostafa"
$ws.Range("G36").Value = "Dr. Marian Samir, Dr. Enas Omran, Dr. Walaa Ghanima"
$ws.Range("G37").Value = "Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida"
$ws.Range("G38").Value = "Dr. Marina Sorial, Dr. Nardine, Dr. Marina Atef, Dr. Neveen Nashaat, Dr. Remon, Dr. Yasmin, Dr. Monica"
$ws.Range("G39").Value = "Dr. Marina Sorial, Dr. Nardine, Dr. Marina Atef, Dr. Neveen Nashaat, Dr. Remon, Dr. Yasmin, Dr. Monica"
$ws.Range("G40").Value = "Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Alshimaa Atef, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Shimaa Ahmad Mekki, Dr. Heba Mahmoud Ali"
$ws.Range("G41").Value = "Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Alshimaa Atef, Dr. Shimaa Ahmad Mekki, Dr. Mohammad El-Tanany, Dr. Hanan Ragab"
$ws.Range("G42").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Eman Tantawi, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad, Dr. Shimaa Ahmad Mekki"
$ws.Range("G43").Value = "Dr. Fatma Elhady, Dr. Nada Mohammad, Dr. Menna tu'Alllah Mohammad, Dr. Abeer Ragab, Dr. Lamiaa Ossama, Dr. Amera Ahmad Saad, Dr. Kerelos Zareef"
$ws.Range("G44").Value = "Dr. Sara Nabil, Dr. Safa Hany"
$ws.Range("G45").Value = "Dr. Amal Awwad, Dr. Safa Hany"
$ws.Range("G47").Value = "Dr. Nourhan Osama, Dr. Amira Ibrahim, Dr. Arwa Al-Sayed, Dr. Maryam Ahmad, Dr. Esraa Mostafa, Dr. Merna Said"
$ws.Range("G48").Value = "Dr. Fatma Shoukry, Dr. Yasmeena Fattoh, Dr. Sarah Abdelmohsen, Dr. Maryam Ahmad, Dr. Eman M. Abo-Sakaya, Dr. Amany Raafat, Dr. Merna Said"
$ws.Range("G49").Value = "Dr. Mariam Toma Gerges, Dr. Mohammad Safwat"
$ws.Range("G56").Value = "Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida"
$ws.Range("G57").Value = "Dr. Marina Sorial, Dr. Nardine, Dr. Marina Atef, Dr. Neveen Nashaat, Dr. Remon, Dr. Yasmin, Dr. Monica"
$ws.Range("G58").Value = "Dr. Marina Sorial, Dr. Nardine, Dr. Marina Atef, Dr. Neveen Nashaat, Dr. Remon, Dr. Yasmin, Dr. Monica"
$ws.Range("G59").Value = "Dr. Amira Sobhy, Dr. Nourhan Mahmoud, Dr. Nesma, Dr. Servinaz Sayed Mohammad, Dr. Mohammad El-Tanany, Dr. Heba Mahmoud Ali, Dr. Asmaa Reda"
$ws.Range("G60").Value = "Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Alshimaa Atef, Dr. Shimaa Ahmad Mekki, Dr. Mohammad El-Tanany, Dr. Hanan Ragab"
$ws.Range("G61").Value = "Dr. Amira Sobhy, Dr. Majorelle Magdy, Dr. Shimaa Ahmad Mekki, Dr. Nahla Nagiub, Dr. Asmaa Reda"
$ws.Range("G63").Value = "Dr. Amal Awwad, Dr. Aya Saeed, Dr. Safa Hany"
$ws.Range("G66").Value = "Dr. Marina Youhanna, Dr. Amira Ibrahim, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Dina Adel, Dr. Eman M. Abo-Sakaya"
$ws.Range("G67").Value = "Dr. Amira Ibrahim, Dr. Yasmeena Fattoh, Dr. Esraa Mostafa"
$ws.Range("G75").Value = "Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida"
$ws.Range("G76").Value = "Dr. Marina Sorial, Dr. Wafaa Ebida, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Eman Samir Gabry"
$ws.Range("G77").Value = "Dr. Marina Sorial, Dr. Nardine, Dr. Marina Atef, Dr. Neveen Nashaat, Dr. Remon, Dr. Yasmin, Dr. Monica"
$ws.Range("G78").Value = "Dr. Amira Sobhy, Dr. Nourhan Mahmoud, Dr. Nesma, Dr. Servinaz Sayed Mohammad, Dr. Mohammad El-Tanany, Dr. Heba Mahmoud Ali, Dr. Asmaa Reda"
$ws.Range("G79").Value = "Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Alshimaa Atef, Dr. Shimaa Ahmad Mekki, Dr. Mohammad El-Tanany, Dr. Hanan Ragab"
$ws.Range("G80").Value = "Dr. Amira Sobhy, Dr. Majorelle Magdy, Dr. Shimaa Ahmad Mekki, Dr. Nahla Nagiub, Dr. Asmaa Reda"
$ws.Range("G81").Value = "Dr. Fatma Elhady, Dr. Nada Mohammad, Dr. Abeer Ragab, Dr. Lamiaa Ossama, Dr. Amera Ahmad Saad"
$ws.Range("G82").Value = "Dr. Yasmin Tarek, Dr. Nourhan Mohammad"
$ws.Range("G83").Value = "Dr. Amal Awwad, Dr. Aya Saeed, Dr. Safa Hany"
$ws.Range("G85").Value = "Dr. Marina Youhanna, Dr. Amira Ibrahim, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Dina Adel, Dr. Eman M. Abo-Sakaya"
$ws.Range("G86").Value = "Dr. Fatma Shoukry, Dr. Yasmeena Fattoh, Dr. Sarah Abdelmohsen, Dr. Maryam Ahmad, Dr. Eman M. Abo-Sakaya, Dr. Amany Raafat, Dr. Merna Said"
$ws.Range("G94").Value = "Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida"
$ws.Range("G95").Value = "Dr. Marina Sorial, Dr. Wafaa Ebida, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Eman Samir Gabry"
$ws.Range("G96").Value = "Dr. Marina Sorial, Dr. Nardine, Dr. Marina Atef, Dr. Neveen Nashaat, Dr. Remon, Dr. Yasmin, Dr. Monica"
$ws.Range("G97").Value = "Dr. Amira Sobhy, Dr. Nourhan Mahmoud, Dr. Nesma, Dr. Servinaz Sayed Mohammad, Dr. Mohammad El-Tanany, Dr. Heba Mahmoud Ali, Dr. Asmaa Reda"
$ws.Range("G98").Value = "Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Alshimaa Atef, Dr. Shimaa Ahmad Mekki, Dr. Mohammad El-Tanany, Dr. Hanan Ragab"
$ws.Range("G99").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Eman Tantawi, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad, Dr. Shimaa Ahmad Mekki"
$ws.Range("G100").Value = "Dr. Fatma Elhady, Dr. Nada Mohammad, Dr. Menna tu'Alllah Mohammad, Dr. Abeer Ragab, Dr. Lamiaa Ossama, Dr. Amera Ahmad Saad, Dr. Kerelos Zareef"
$ws.Range("G101").Value = "Dr. Amal Awwad, Dr. Aya Saeed, Dr. Safa Hany"
$ws.Range("G102").Value = "Dr. Amal Awwad, Dr. Safa Hany"
$ws.Range("G104").Value = "Dr. Nourhan Osama, Dr. Amira Ibrahim, Dr. Arwa Al-Sayed, Dr. Maryam Ahmad, Dr. Esraa Mostafa, Dr. Merna Said"
$ws.Range("G105").Value = "Dr. Amira Ibrahim, Dr. Basma Hamed"
$ws.Range("G113").Value = "Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida"
$ws.Range("G115").Value = "Dr. Marina Sorial, Dr. Nardine, Dr. Marina Atef, Dr. Neveen Nashaat, Dr. Remon, Dr. Yasmin, Dr. Monica"
